$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 14.66769855181886
